# Auto-generated edit script applying the row-content update for Översikt HELSINGBORG
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Cells.Item(2,1).Value2 = "A 45325-2025"
$ws.Cells.Item(2,2).Value2 = 45922
$ws.Cells.Item(2,3).Value2 = 46066
$ws.Cells.Item(2,7).Value2 = 1.6
$ws.Cells.Item(2,8).Value2 = 4
$ws.Cells.Item(2,9).Value2 = 0
$ws.Cells.Item(2,10).Value2 = 1
$ws.Cells.Item(2,15).Value2 = 1
$ws.Cells.Item(2,18).Value2 = "Nordlig buksimmare`r`nStörre vattensalamander`r`nÅkergroda`r`nMindre vattensalamander`r`nVanlig groda"
$ws.Cells.Item(2,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/artfynd/A 45325-2025 artfynd.xlsx", "A 45325-2025")'
$ws.Cells.Item(2,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/kartor/A 45325-2025 karta.png", "A 45325-2025")'
$ws.Cells.Item(2,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomål/A 45325-2025 FSC-klagomål.docx", "A 45325-2025")'
$ws.Cells.Item(2,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomålsmail/A 45325-2025 FSC-klagomål mail.docx", "A 45325-2025")'
$ws.Cells.Item(2,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsyn/A 45325-2025 tillsynsbegäran.docx", "A 45325-2025")'
$ws.Cells.Item(2,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsynsmail/A 45325-2025 tillsynsbegäran mail.docx", "A 45325-2025")'
$ws.Cells.Item(2,26).ClearContents()
$ws.Rows.Item(2).RowHeight = 15
# --- Row 3 ---
$ws.Cells.Item(3,1).Value2 = "A 13467-2023"
$ws.Cells.Item(3,2).Value2 = 45005
$ws.Cells.Item(3,3).Value2 = 46066
$ws.Cells.Item(3,7).Value2 = 2.3
$ws.Cells.Item(3,8).Value2 = 2
$ws.Cells.Item(3,9).Value2 = 2
$ws.Cells.Item(3,10).Value2 = 3
$ws.Cells.Item(3,15).Value2 = 3
$ws.Cells.Item(3,18).Value2 = "Gulsparv`r`nHypoxylon petriniae`r`nKråka`r`nGrå skärelav`r`nGulnål"
$ws.Cells.Item(3,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/artfynd/A 13467-2023 artfynd.xlsx", "A 13467-2023")'
$ws.Cells.Item(3,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/kartor/A 13467-2023 karta.png", "A 13467-2023")'
$ws.Cells.Item(3,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomål/A 13467-2023 FSC-klagomål.docx", "A 13467-2023")'
$ws.Cells.Item(3,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomålsmail/A 13467-2023 FSC-klagomål mail.docx", "A 13467-2023")'
$ws.Cells.Item(3,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsyn/A 13467-2023 tillsynsbegäran.docx", "A 13467-2023")'
$ws.Cells.Item(3,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsynsmail/A 13467-2023 tillsynsbegäran mail.docx", "A 13467-2023")'
$ws.Cells.Item(3,26).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/fåglar/A 13467-2023 prioriterade fågelarter.docx", "A 13467-2023")'
$ws.Rows.Item(3).RowHeight = 15
# --- Row 4 ---
$ws.Cells.Item(4,3).Value2 = 46066
# --- Row 5 ---
$ws.Cells.Item(5,3).Value2 = 46066
# --- Row 6 ---
$ws.Cells.Item(6,1).Value2 = "A 2593-2024"
$ws.Cells.Item(6,2).Value2 = 45313.69204861111
$ws.Cells.Item(6,3).Value2 = 46066
$ws.Cells.Item(6,7).Value2 = 2.3
# --- Row 7 ---
$ws.Cells.Item(7,1).Value2 = "A 7827-2026"
$ws.Cells.Item(7,2).Value2 = 46062.63958333333
$ws.Cells.Item(7,3).Value2 = 46066
$ws.Cells.Item(7,7).Value2 = 2.1
# --- Row 8 ---
$ws.Cells.Item(8,1).Value2 = "A 7814-2026"
$ws.Cells.Item(8,2).Value2 = 46062.61388888889
$ws.Cells.Item(8,3).Value2 = 46066
$ws.Cells.Item(8,7).Value2 = 1.1
# --- Row 9 ---
$ws.Cells.Item(9,1).Value2 = "A 12651-2022"
$ws.Cells.Item(9,2).Value2 = 44641
$ws.Cells.Item(9,3).Value2 = 46066
$ws.Cells.Item(9,7).Value2 = 3.2
# --- Row 10 ---
$ws.Cells.Item(10,1).Value2 = "A 5792-2024"
$ws.Cells.Item(10,2).Value2 = 45335
$ws.Cells.Item(10,3).Value2 = 46066
$ws.Cells.Item(10,7).Value2 = 5.6
# --- Row 11 ---
$ws.Cells.Item(11,1).Value2 = "A 13651-2023"
$ws.Cells.Item(11,2).Value2 = 45006
$ws.Cells.Item(11,3).Value2 = 46066
$ws.Cells.Item(11,7).Value2 = 2.2
# --- Row 12 ---
$ws.Cells.Item(12,1).Value2 = "A 8194-2025"
$ws.Cells.Item(12,2).Value2 = 45708
$ws.Cells.Item(12,3).Value2 = 46066
$ws.Cells.Item(12,7).Value2 = 1.9
# --- Row 13 ---
$ws.Cells.Item(13,1).Value2 = "A 35642-2023"
$ws.Cells.Item(13,2).Value2 = 45147
$ws.Cells.Item(13,3).Value2 = 46066
$ws.Cells.Item(13,7).Value2 = 1.2
# --- Row 14 ---
$ws.Cells.Item(14,3).Value2 = 46066
# --- Row 15 ---
$ws.Cells.Item(15,1).Value2 = "A 28288-2023"
$ws.Cells.Item(15,2).Value2 = 45099.6349537037
$ws.Cells.Item(15,3).Value2 = 46066
$ws.Cells.Item(15,7).Value2 = 0.5
# --- Row 16 ---
$ws.Cells.Item(16,1).Value2 = "A 7333-2025"
$ws.Cells.Item(16,2).Value2 = 45703.35899305555
$ws.Cells.Item(16,3).Value2 = 46066
$ws.Cells.Item(16,7).Value2 = 0.9
